$d = $word.ActiveDocument

# Hunk 1: " on July 01, 2022." -> " on July 02, 2022."
$d.Content.Find.Execute(" on July 01, 2022.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " on July 02, 2022.", 1)

# Hunk 2: standalone bold "July 01, 2022" -> "July 02, 2022"
$d.Content.Find.Execute("July 01, 2022", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "July 02, 2022", 1)

# Hunk 3: "August 30, 2022" -> "August 31, 2022"
$d.Content.Find.Execute("August 30, 2022", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "August 31, 2022", 1)

# Hunk 4: " license is suspended from July 01, 2022" -> " license is suspended from July 02, 2022"
$d.Content.Find.Execute(" license is suspended from July 01, 2022", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " license is suspended from July 02, 2022", 1)
